$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "87.819.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.177.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.94%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.39%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "610.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.48%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.387"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.47%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.674"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.66%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.173.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -9.68%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -9.03%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.760.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.83%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.691.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -7.55%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.185.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.23%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.70%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.93%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "414.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.12%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -8.80%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -6.83%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.343.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.64%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000133"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "73.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.06%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -13.55%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "546.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -9.51%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -9.84%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.86%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "21.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.130"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.50%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "21.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.42%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -8.07%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "174.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.72%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.86%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -10.08%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -8.14%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.66%  "
$ws.Range("E51").Style = "Normal"
